# Update "Forecast Comparison" sheet with the correct forecast output:
# insert a new "Week_Start_Date" column after "Week", and normalize the
# week labels (drop the zero-padding, e.g. "W01" -> "W1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column before the current column B (ASIN), shifting
# ASIN..is_holiday_week from B:I to C:J.
$ws.Columns.Item(2).Insert()

# Header for the freshly inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# The new column holds plain-text ISO dates (one per week, starting
# 2025-01-05) - force text formatting so Excel doesn't coerce the
# assigned strings into date serial numbers.
$ws.Range("B2:B17").NumberFormat = "@"

$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Drop the leading zero in the week label ("W01" -> "W1", "W09" -> "W9";
    # "W10".."W16" are already unpadded and stay the same).
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)

    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
}

Write-Host "Inserted Week_Start_Date column and refreshed week labels."
